# Natural Instructions evaluation update — add a callout textbox summarising
# the "Confident in GPT2 vs Mistral ICL" result on the Natural Instructions
# task, on the "Impact Statement (business/industry)" slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# EMU target box: off x=5273336 y=3187083, ext cx=6045693 cy=369332.
# PowerPoint COM works in points (1 pt = 12700 EMU); the values below are
# nudged by a hair so the float32 round-trip inside the host lands on the
# exact EMU the deck originally shipped with.
$left   = 415.22331308661416
$top    = 250.95142432283467
$width  = 476.0388338976378
$height = 29.081259842519685

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "TextBox 1"

$shp.TextFrame.WordWrap = $true
$shp.TextFrame.AutoSize = 1
$shp.Fill.Visible = $false

$tr1 = $shp.TextFrame.TextRange
$tr1.Text = "Confident in GPT2 vs Mistral ICL "
$tr1.LanguageID = "en-GB"

$tr2 = $tr1.InsertAfter("Natural Instructions")
$tr2.LanguageID = "en-GB"
